$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.615.79"
$ws.Range("E2").Value = "  +1.24%  "

$ws.Range("D3").Value = "3.789.86"
$ws.Range("E3").Value = "  +1.95%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "420.78"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.35%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "128.49"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.14%  "

$ws.Range("D7").Value = "3.787.33"
$ws.Range("E7").Value = "  +2.09%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.602"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.32%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.12%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.718"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.79%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.160"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.57%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000346"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.36%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "39.99"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.81%  "

$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "4.395.46"
$ws.Range("E14").Value = "  +1.67%  "

$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "10.04"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.10%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.68"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +20.22%  "

$ws.Range("E17").Value = "  -0.93%  "

$ws.Range("D18").Value = "3.795.88"
$ws.Range("E18").Value = "  +1.69%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "19.44"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.52%  "

$ws.Range("D20").Value = "66.899.45"
$ws.Range("E20").Value = "  +1.04%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.08"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.18%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "403.24"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.73%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.18"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.10%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.48"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.23%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.99"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.37%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "36.81"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.98%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "5.55"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +7.25%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.20"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.53%  "

$ws.Range("E29").Value = "  -1.65%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "720.57"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.91%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.50"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +20.87%  "

$ws.Range("B32").Value = "Toncoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.76"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.99%  "

$ws.Range("B33").Value = "Cosmos"
$ws.Range("C33").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "12.32"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.67%  "

$ws.Range("E34").Value = "  +0.08%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.01%  "

$ws.Range("E36").Value = "  -6.28%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "38.33"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -7.63%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "54.87"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.96%  "

$ws.Range("D39").Value = "0.0₃0750"
$ws.Range("E39").Value = "  +10.15%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.93"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +16.07%  "

$ws.Range("E41").Value = "  -4.56%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.91"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.31%  "

$ws.Range("E43").Value = "  +0.29%  "

$ws.Range("E44").Value = "  -5.11%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.32"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.89%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "143.46"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.29%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.10"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.89%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.02"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.47%  "

$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.77"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.80%  "

$ws.Range("B50").Value = "WEMIXToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.55"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.47%  "

$ws.Range("B51").Value = "TheGraph"
$ws.Range("C51").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.305"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.39%  "
